$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "network"
$ws.Range("C1").Value = "site-id"
$ws.Range("D1").Value = "site-name"

# Column A (ids) - numeric, unchanged values but set anyway for completeness
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2

# Column B (network) - write column-wise to match shared string ordering
$ws.Range("B2").Value = "10.10.10.0/24"
$ws.Range("B3").Value = "10.10.11.0/24"

# Column C (site-id)
$ws.Range("C2").Value = "abcf02e3-8c11-4aef-b3f7-7a5284471c4f"
$ws.Range("C3").Value = "abc74fc3-8819-4cf0-b318-b7be37b21b7d"

# Column D (site-name)
$ws.Range("D2").Value = "Site 1"
$ws.Range("D3").Value = "Site 2"

# Update the active selection to A2 (was A3)
$ws.Range("A2").Select()
